$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) number format on the D/E cells we are about to
# rewrite, so numeric-looking values (e.g. "244.92", "1.001") are stored
# as text rather than being auto-coerced to numbers by Excel -- matching
# the original inline-string cell type. ClearFormats() afterwards removes
# the temporary formatting so cell styling stays untouched.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.571.47'
$ws.Range('E2').Value = '  -0.14%  '

$ws.Range('D3').Value = '1.914.89'
$ws.Range('E3').Value = '  -0.46%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '244.92'
$ws.Range('E5').Value = '  -0.82%  '

$ws.Range('E6').Value = '  +0.00%  '

$ws.Range('D7').Value = '0.4842'
$ws.Range('E7').Value = '  +1.97%  '

$ws.Range('D8').Value = '0.2890'

$ws.Range('D9').Value = '0.06794'
$ws.Range('E9').Value = '  -0.43%  '

$ws.Range('D10').Value = '111.74'
$ws.Range('E10').Value = '  +6.25%  '

$ws.Range('D11').Value = '19.36'
$ws.Range('E11').Value = '  +5.27%  '

$ws.Range('D12').Value = '1.913.87'
$ws.Range('E12').Value = '  -0.52%  '

$ws.Range('D13').Value = '0.07579'
$ws.Range('E13').Value = '  -1.37%  '

$ws.Range('D14').Value = '5.404'
$ws.Range('E14').Value = '  +1.41%  '

$ws.Range('D15').Value = '0.6703'
$ws.Range('E15').Value = '  +0.29%  '

$ws.Range('D16').Value = '293.85'
$ws.Range('E16').Value = '  +1.46%  '

$ws.Range('D17').Value = '30.553.01'
$ws.Range('E17').Value = '  -0.25%  '

$ws.Range('E18').Value = '  +0.47%  '

$ws.Range('D19').Value = '0.000007596'
$ws.Range('E19').Value = '  -0.20%  '

$ws.Range('D21').Value = '5.512'
$ws.Range('E21').Value = '  -1.29%  '

$ws.Range('D22').Value = '2.161.51'
$ws.Range('E22').Value = '  -0.56%  '

$ws.Range('D23').Value = '1.001'

$ws.Range('D24').Value = '6.416'
$ws.Range('E24').Value = '  -0.45%  '

$ws.Range('D25').Value = '9.462'
$ws.Range('E25').Value = '  -0.38%  '

$ws.Range('D26').Value = '166.06'
$ws.Range('E26').Value = '  -0.43%  '

$ws.Range('E27').Value = '  -4.32%  '

$ws.Range('E28').Value = '  -1.59%  '

$ws.Range('E29').Value = '  -0.74%  '

$ws.Range('D30').Value = '1.441'
$ws.Range('E30').Value = '  +2.74%  '

$ws.Range('E31').Value = '  -1.29%  '

$ws.Range('D32').Value = '4.046'
$ws.Range('E32').Value = '  +0.10%  '

$ws.Range('D33').Value = '0.04977'
$ws.Range('E33').Value = '  -1.06%  '

$ws.Range('D34').Value = '0.7338'
$ws.Range('E34').Value = '  +0.52%  '

$ws.Range('D35').Value = '1.140'
$ws.Range('E35').Value = '  -0.22%  '

$ws.Range('D36').Value = '2.719'
$ws.Range('E36').Value = '  -0.67%  '

$ws.Range('D37').Value = '0.02026'
$ws.Range('E37').Value = '  -1.89%  '

$ws.Range('D38').Value = '2.682'
$ws.Range('E38').Value = '  -0.30%  '

$ws.Range('E39').Value = '  -1.45%  '

$ws.Range('D40').Value = '109.33'
$ws.Range('E40').Value = '  -2.13%  '

$ws.Range('D41').Value = '0.4428'
$ws.Range('E41').Value = '  +0.87%  '

$ws.Range('E42').Value = '  -0.94%  '

$ws.Range('D43').Value = '5.830'
$ws.Range('E43').Value = '  -1.58%  '

$ws.Range('E44').Value = '  +0.02%  '

$ws.Range('D45').Value = '69.30'
$ws.Range('E45').Value = '  +2.42%  '

$ws.Range('D46').Value = '7.204'
$ws.Range('E46').Value = '  -1.09%  '

$ws.Range('D47').Value = '48.58'
$ws.Range('E47').Value = '  -0.10%  '

$ws.Range('D48').Value = '9.249'
$ws.Range('E48').Value = '  -0.85%  '

$ws.Range('D49').Value = '0.1225'
$ws.Range('E49').Value = '  -1.52%  '

$ws.Range('D50').Value = '0.2500'
$ws.Range('E50').Value = '  +0.48%  '

$ws.Range('D51').Value = '34.74'
$ws.Range('E51').Value = '  -0.57%  '

$ws.Range("D2:E51").ClearFormats()
